$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: extend the notes text and bump hours spent from 1 to 2
$ws.Range("D11").Value = "Defined the first draft of the data model for reporting work package and started building the properties file handling."
$ws.Range("C11").Value = 2
$ws.Rows.Item(11).RowHeight = 28.5

# New row 12: next day's entry
$ws.Range("A12").Value = 42902
$ws.Range("B12").Formula = "=A12"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "fleshed out the properties handling."

# Copy the date/formula column formatting from row 11 down into row 12
# (values/formulas already set above; this only brings over the number formats)
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)

# Move the active selection, matching the saved workbook state
$ws.Range("D14").Select()
